# Verify ENTER ACCOUNT INFORMATION header is visible test
# -----------------------------------------------------------------
# This script reproduces, via Excel COM-interop, the edits that were
# made to dynamicdata.xlsx:
#   - DynamicData sheet: fill in rows 3 & 4 (email / password pairs for
#     TestUser2 / TestUser3) and turn the e-mail cells into mailto
#     hyperlinks, matching the existing A1/A2 pattern.
#   - Sheet1: turn the two e-mail columns for TestUser2 / TestUser3
#     rows into mailto hyperlinks as well (Hyperlink style).
#   - Update the saved sheet views (zoom / selection) on both sheets.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- DynamicData sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("DynamicData")

# Row 4 first (TestUser3), then row 3 (TestUser2) - matches the order the
# hyperlink relationships were created in the authored workbook.
$ws1.Range("A4").Value = "user3_@example.com"
$ws1.Range("B4").Value = "MyPass789"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "mailto:user3_@example.com")
$ws1.Range("A4").Style = "Hyperlink"
$ws1.Range("B4").Style = "Normal"

$ws1.Range("A3").Value = "user2_@example.com"
$ws1.Range("B3").Value = "Secret456"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:user2_@example.com")
$ws1.Range("A3").Style = "Hyperlink"
$ws1.Range("B3").Style = "Normal"

# Updated view: zoom 106%, selection on A7
[void]$ws1.Range("A7").Select()
$excel.ActiveWindow.Zoom = 106

# ---- Sheet1 (tabular data) ---------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet1")
[void]$ws2.Activate()

# The e-mail values are already present in B3/B4; just wire up mailto
# hyperlinks and apply the Hyperlink style, same as Excel would do when
# you use Insert > Link on a cell that already contains text.
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:user2_{{unique}}@example.com")
$ws2.Range("B3").Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Range("B4"), "mailto:user3_{{unique}}@example.com")
$ws2.Range("B4").Style = "Hyperlink"

# Updated view: selection moves to B4 (also drops the old topLeftCell scroll)
[void]$ws2.Range("B4").Select()
